$d = $word.ActiveDocument

# Merge "Use" + proofErr marks + " Case" into a single run reading
# "nach Use Case" (the search text starts inside the preceding run so the
# whole span, including the spell-check proofErr markers around "Use", is
# rebuilt as one clean run).
$d.Content.Find.Execute("ch Use Case", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "ch nach Use Case", 2)

# Remove the now-redundant "des Systems nach " lead-in text.
$d.Content.Find.Execute("des Systems nach ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2)

# Re-locate the cell text and split the run right after "Zustand " by
# dropping the _GoBack bookmark there - this moves the bookmark from the
# trailing empty paragraph (after the table) to between the two runs.
$t = $d.Tables.Item(1)
$cell = $t.Cell(8, 1)
$splitAt = $cell.Range.Start + "Zustand ".Length
$r = $d.Range($splitAt, $splitAt)
$d.Bookmarks.Add("_GoBack", $r)
